# Prueba a cuentas bancarias 0.1
# Adds an "email" row (row 6) with hyperlinked e-mail addresses for the
# first three users in the sheet, tweaks a couple of column widths and
# moves the active selection from D8 to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row 6 : email header + 3 mailto hyperlinks ------------------
$ws.Range("A6").Value = "email"
$ws.Range("B6").Value = "lzapata@edeq.com"
$ws.Range("C6").Value = "mgarcia@edeq.com"
$ws.Range("D6").Value = "dflores@edeq.com"

# D6 keeps the sheet's "text" number format, right aligned, same as the
# other numeric/right-aligned columns on the sheet - set this before the
# hyperlink style is applied so the resulting xf stays linked to the
# Hyperlink cell style (xfId=1).
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").HorizontalAlignment = -4152

# Turn B6/C6/D6 into real hyperlinks (mailto:) - this is what introduces
# the new "Hyperlink" font/style into styles.xml.
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:lzapata@edeq.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:mgarcia@edeq.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:dflores@edeq.com") | Out-Null

# --- column widths -----------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 17.333333333333332
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668

# --- selection moves from D8 to D7 -------------------------------------
$ws.Range("D7").Select() | Out-Null
